$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тест-кейсы на автоматизацию")

# --- Row 14: tweak the existing "Create 2 Test Cases..." test case row ---
$ws.Range("B14").Value = "Create 2 Test Cases and add to Test Plan as Suite"
$ws.Range("C14").Value = "E2E"

# --- Row 15: new test case row, added as a Suite-style sibling row ---
# Set B15 before D14/D15 so new shared-string entries come out in the same
# order as the target workbook ("Create 3 Test Cases..." before "To do").
$ws.Range("B15").Value = "Create 3 Test Cases in 1 Suite and add 2 to Test Plan"
$ws.Range("B15").Font.Bold = $false
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("B15").VerticalAlignment = -4108
$ws.Range("B15").WrapText = $true
$ws.Range("C15").Value = "E2E"

# --- Status column ("Статус") for both new rows ---
$ws.Range("D14").Value = "To do"
$ws.Range("D15").Value = "To do"

# --- Drop leftover stray formatting on D17 (no value, no longer needed) ---
$ws.Range("D17").Clear()

# --- Restore selection cursor where the author left it ---
$ws.Range("D16").Select()
